$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column at G ("Collection"), shifting Length..Universal Fish
#    ID (old G:AC) one column to the right (new H:AD). Excel's own Insert()
#    already takes care of: dimension, col widths/shift, row spans, cell
#    values/styles, pane/selection topLeftCell.
# ---------------------------------------------------------------------------
$ws.Columns("G").Insert()

# New header cell/value for the inserted column.
$ws.Range("G1").Value = "Collection"

# ---------------------------------------------------------------------------
# 2. Cell comments do not move with Insert(), so shift them manually one
#    column to the right, from the rightmost one down to G, so we never
#    clobber a comment we still need to read.
# ---------------------------------------------------------------------------
$pairs = @(
    @("W","X"), @("V","W"), @("T","U"), @("S","T"), @("R","S"), @("Q","R"),
    @("N","O"), @("M","N"), @("L","M"), @("K","L"), @("J","K"), @("I","J"), @("H","I"), @("G","H")
)

foreach ($p in $pairs) {
    $srcCol = $p[0]
    $dstCol = $p[1]
    $srcComment = $ws.Range($srcCol + "1").Comment
    if ($srcComment -ne $null) {
        $txt = $srcComment.Text()
        $dstComment = $ws.Range($dstCol + "1").Comment
        if ($dstComment -ne $null) {
            $dstComment.Text($txt)
        } else {
            $ws.Range($dstCol + "1").AddComment($txt)
        }
    }
}

# G1 still holds the stale "Length" comment (its text was already copied to
# H1 above) -- overwrite it with the brand new note for the Collection column.
$g1Comment = $ws.Range("G1").Comment
if ($g1Comment -ne $null) {
    $g1Comment.Text("Eg. 1999 FP")
} else {
    $ws.Range("G1").AddComment("Eg. 1999 FP")
}

# Q1 and V1 are left-over sources that nothing shifted into -- drop them.
$q1Comment = $ws.Range("Q1").Comment
if ($q1Comment -ne $null) { $q1Comment.Delete() }
$v1Comment = $ws.Range("V1").Comment
if ($v1Comment -ne $null) { $v1Comment.Delete() }

# ---------------------------------------------------------------------------
# 3. Conditional formatting on the "Cond.Fact." column tracked I1; move its
#    range to J1 (new position) while keeping the same dxf/priority.
# ---------------------------------------------------------------------------
$condFormats = $ws.Range("I1").FormatConditions
for ($i = 1; $i -le $condFormats.Count; $i++) {
    $condFormats.Item($i).ModifyAppliesToRange($ws.Range("J1"))
}

# ---------------------------------------------------------------------------
# 4. AutoFilter range grows from A1:AC1 to A1:AD1. Toggle off/on so the
#    stored range is refreshed instead of just being removed.
# ---------------------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:AD1").AutoFilter()

# ---------------------------------------------------------------------------
# 5. The hidden defined name backing the autofilter needs the same update.
# ---------------------------------------------------------------------------
$filterName = $wb.Names.Item(1)
$filterName.RefersTo = "=Template!`$A`$1:`$AD`$1"

# ---------------------------------------------------------------------------
# 6. Restore the active selection to the newly inserted column's header.
# ---------------------------------------------------------------------------
$ws.Range("G1").Select()
